$wb = $excel.ActiveWorkbook

function Set-TextCell($ws, $row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $text
    $c.Style = "Normal"
}

function Set-EmptyCell($ws, $row, $col) {
    $ws.Cells.Item($row, $col).Style = "Normal"
}

$ws = $wb.Worksheets.Item("cases")

# Header: new column Y (col 25), continuing the date sequence
Set-TextCell $ws 1 25 "2020-04-26"

# Rows 2-21: empty placeholder cells in new column Y
for ($r = 2; $r -le 21; $r++) {
    Set-EmptyCell $ws $r 25
}

# Row 22: Y22 empty placeholder
Set-EmptyCell $ws 22 25

# Row 23: B23 newly populated, Y23 empty placeholder
$ws.Cells.Item(23, 2).Value = 61888
Set-EmptyCell $ws 23 25

# Rows 24-36: new forecast values in column Y
$ws.Cells.Item(24, 25).Value = 67049
$ws.Cells.Item(25, 25).Value = 71646
$ws.Cells.Item(26, 25).Value = 77345
$ws.Cells.Item(27, 25).Value = 85083
$ws.Cells.Item(28, 25).Value = 91409
$ws.Cells.Item(29, 25).Value = 99959
$ws.Cells.Item(30, 25).Value = 105199
$ws.Cells.Item(31, 25).Value = 115887
$ws.Cells.Item(32, 25).Value = 123388
$ws.Cells.Item(33, 25).Value = 130387
$ws.Cells.Item(34, 25).Value = 137236
$ws.Cells.Item(35, 25).Value = 143194
$ws.Cells.Item(36, 25).Value = 147517

# Row 37: brand new row (new forecast date)
Set-TextCell $ws 37 1 "2020-05-10"
for ($c = 2; $c -le 24; $c++) {
    Set-EmptyCell $ws 37 $c
}
$ws.Cells.Item(37, 25).Value = 153226

$ws = $wb.Worksheets.Item("deaths")

# Header: new column Y (col 25), continuing the date sequence
Set-TextCell $ws 1 25 "2020-04-26"

# Rows 2-21: empty placeholder cells in new column Y
for ($r = 2; $r -le 21; $r++) {
    Set-EmptyCell $ws $r 25
}

# Row 22: Y22 empty placeholder
Set-EmptyCell $ws 22 25

# Row 23: B23 newly populated, Y23 empty placeholder
$ws.Cells.Item(23, 2).Value = 4205
Set-EmptyCell $ws 23 25

# Rows 24-36: new forecast values in column Y
$ws.Cells.Item(24, 25).Value = 4761
$ws.Cells.Item(25, 25).Value = 5095
$ws.Cells.Item(26, 25).Value = 5604
$ws.Cells.Item(27, 25).Value = 6253
$ws.Cells.Item(28, 25).Value = 6822
$ws.Cells.Item(29, 25).Value = 7505
$ws.Cells.Item(30, 25).Value = 8064
$ws.Cells.Item(31, 25).Value = 8989
$ws.Cells.Item(32, 25).Value = 9741
$ws.Cells.Item(33, 25).Value = 10498
$ws.Cells.Item(34, 25).Value = 11176
$ws.Cells.Item(35, 25).Value = 11772
$ws.Cells.Item(36, 25).Value = 12142

# Row 37: brand new row (new forecast date)
Set-TextCell $ws 37 1 "2020-05-10"
for ($c = 2; $c -le 24; $c++) {
    Set-EmptyCell $ws 37 $c
}
$ws.Cells.Item(37, 25).Value = 12779
